$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Add the new ToDo item as row 24
$ws.Range("A24").Value = "Standardize on confirmation dialogs"
$ws.Range("B24").Value = "There are many places in the system where there are ""remove"" or ""delete"" buttons. Some of these buttons will open a dialog box asking the user to confirm the action, but others do not. It should be the same behaviour everywhere."
$ws.Range("C24").Value = "OPEN"

# Match the row height used by other multi-line description rows (45pt)
$ws.Range("A24:C24").RowHeight = 45

# Scroll the view down a bit and select the new status cell, as in the source workbook
$excel.ActiveWindow.ScrollRow = 9
$ws.Range("C24").Select()
